# Revert "Powerpoint writer: consolidate text run nodes."
#
# The writer used to merge a trailing space onto the end of the
# preceding word's run, e.g. a single run "Slide " immediately
# followed by a run "1 ". That consolidation is undone here: every
# run whose text ends with a single trailing space must have that
# space carved out into its own, separate run, e.g.
#   "Slide " + "1 " + "(Content)"   ->   "Slide" + " " + "1" + " " + "(Content)"
#
# This host's TextRange has no working Runs() collection to walk, but
# re-assigning a TextRange.Characters(start,length) sub-range's own
# text back onto itself forces the host to re-partition the backing
# runs at that exact character boundary - exactly like real
# PowerPoint does when you select text and retype the same thing over
# it. So for every "word " that used to be its own run (identified by
# its text below) we look up where it sits in the shape's current
# text and re-assign everything up to - but not including - its
# trailing space, which pushes that space out into a run of its own.
#
# NB: this host's PowerShell only accepts *positional* function
# arguments - named arguments (`-Foo bar`) are silently swallowed -
# so the helper below is always invoked positionally.

function Split-TrailingSpace($TextRange, $WordWithSpace) {
    $text = $TextRange.Text
    $at = $text.IndexOf($WordWithSpace)
    if ($at -lt 0) {
        throw "Split-TrailingSpace: '$WordWithSpace' not found in '$text'"
    }
    $start = $at + 1                       # 1-based start of the word
    $length = $WordWithSpace.Length - 1    # word length, excluding its trailing space
    $word = $TextRange.Characters($start, $length)
    $word.Text = $word.Text
}

$p = $ppt.ActivePresentation

# Each entry: slide index, shape index, and the runs (identified by
# their original, space-terminated text) that need splitting, in
# left-to-right order.
$targets = @(
    @{ Slide = 1;  Shape = 1; Words = @("Slide ", "1 ") },
    @{ Slide = 2;  Shape = 1; Words = @("Slide ", "2 ") },
    @{ Slide = 3;  Shape = 1; Words = @("Slide ", "3 ") },
    @{ Slide = 4;  Shape = 1; Words = @("Slide ", "4 ") },
    @{ Slide = 5;  Shape = 1; Words = @("Slide ", "5 ", "(Two ") },
    @{ Slide = 6;  Shape = 1; Words = @("Slide ", "6 ", "(Two ", "Content ") },
    @{ Slide = 6;  Shape = 3; Words = @("an ") },
    @{ Slide = 7;  Shape = 1; Words = @("Slide ", "7 ", "(Content ", "with ") },
    @{ Slide = 7;  Shape = 4; Words = @("An ") },
    @{ Slide = 8;  Shape = 1; Words = @("Slide ", "8 ") },
    @{ Slide = 8;  Shape = 4; Words = @("An ") },
    @{ Slide = 9;  Shape = 1; Words = @("Slide ", "10 ") },
    @{ Slide = 10; Shape = 1; Words = @("Slide ", "11 ") },
    @{ Slide = 11; Shape = 1; Words = @("Slide ", "12 ") }
)

foreach ($target in $targets) {
    $slide = $p.Slides.Item($target.Slide)
    $shape = $slide.Shapes.Item($target.Shape)
    $tr = $shape.TextFrame.TextRange
    foreach ($word in $target.Words) {
        Split-TrailingSpace $tr $word
    }
}
